$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 'فصل سوم منتهی به 1399/09'
$ws.Range("E8").Value = 'فصل چهارم منتهی به 1399/12'
$ws.Range("F8").Value = 'فصل اول منتهی به 1400/03'
$ws.Range("G8").Value = 'فصل دوم منتهی به 1400/06'
$ws.Range("H8").Value = 'فصل سوم منتهی به 1400/09'
$ws.Range("I8").Value = 'فصل چهارم منتهی به 1400/12'
$ws.Range("J8").Value = 'فصل اول منتهی به 1401/03'
$ws.Range("K8").Value = 'فصل دوم منتهی به 1401/06'
$ws.Range("L8").Value = 'فصل سوم منتهی به 1401/09'
$ws.Range("M8").Value = 'فصل چهارم منتهی به 1401/12'
$ws.Range("D9").Value = '1399-10-30'
$ws.Range("E9").Value = '1401-02-27 (11)'
$ws.Range("F9").Value = '1400-05-05 (2)'
$ws.Range("G9").Value = '1400-09-29 (3)'
$ws.Range("H9").Value = '1400-11-06 (2)'
$ws.Range("I9").Value = '1401-10-30 (6)'
$ws.Range("J9").Value = '1401-04-29'
$ws.Range("K9").Value = '1401-09-15 (2)'
$ws.Range("L9").Value = '1401-10-30'
$ws.Range("M9").Value = '1402-02-28'
$ws.Range("D12").Value = 116611
$ws.Range("E12").Value = 99419
$ws.Range("F12").Value = 82581
$ws.Range("G12").Value = 158239
$ws.Range("H12").Value = 179404
$ws.Range("I12").Value = 52176
$ws.Range("J12").Value = 48638
$ws.Range("K12").Value = 479672
$ws.Range("L12").Value = 144957
$ws.Range("M12").Value = 664708
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 4460942
$ws.Range("E14").Value = 5229519
$ws.Range("F14").Value = 4857565
$ws.Range("G14").Value = 5575229
$ws.Range("H14").Value = 5640104
$ws.Range("I14").Value = 6677359
$ws.Range("J14").Value = 6517472
$ws.Range("K14").Value = 7256526
$ws.Range("L14").Value = 8620077
$ws.Range("M14").Value = 10423419
$ws.Range("D15").Value = 779013
$ws.Range("E15").Value = 1096698
$ws.Range("F15").Value = 531171
$ws.Range("G15").Value = 1526709
$ws.Range("H15").Value = 1212226
$ws.Range("I15").Value = 2090105
$ws.Range("J15").Value = 1621280
$ws.Range("K15").Value = 1699303
$ws.Range("L15").Value = 1581702
$ws.Range("M15").Value = 2144181
$ws.Range("D16").Value = 736892
$ws.Range("E16").Value = 808867
$ws.Range("F16").Value = 1233143
$ws.Range("G16").Value = 878240
$ws.Range("H16").Value = 1358719
$ws.Range("I16").Value = 600485
$ws.Range("J16").Value = 1242405
$ws.Range("K16").Value = 764853
$ws.Range("L16").Value = 627329
$ws.Range("M16").Value = 1432913
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("D18").Value = 6093458
$ws.Range("E18").Value = 7234503
$ws.Range("F18").Value = 6704460
$ws.Range("G18").Value = 8138417
$ws.Range("H18").Value = 8390453
$ws.Range("I18").Value = 9420125
$ws.Range("J18").Value = 9429795
$ws.Range("K18").Value = 10200354
$ws.Range("L18").Value = 10974065
$ws.Range("M18").Value = 14665221
$ws.Range("D19").Value = 3418
$ws.Range("E19").Value = 8215
$ws.Range("F19").Value = 8215
$ws.Range("G19").Value = 4510
$ws.Range("H19").Value = 3027
$ws.Range("I19").Value = 5308
$ws.Range("J19").Value = 4370
$ws.Range("K19").Value = 5134
$ws.Range("L19").Value = 3503
$ws.Range("M19").Value = 5824
$ws.Range("D20").Value = 607269
$ws.Range("E20").Value = 607270
$ws.Range("F20").Value = 628781
$ws.Range("G20").Value = 886461
$ws.Range("H20").Value = 907170
$ws.Range("I20").Value = 907171
$ws.Range("J20").Value = 902041
$ws.Range("K20").Value = 886731
$ws.Range("L20").Value = 952258
$ws.Range("M20").Value = 941975
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("D22").Value = 380418
$ws.Range("E22").Value = 452852
$ws.Range("F22").Value = 486421
$ws.Range("G22").Value = 556371
$ws.Range("H22").Value = 1552443
$ws.Range("I22").Value = 1703083
$ws.Range("J22").Value = 1714222
$ws.Range("K22").Value = 1848872
$ws.Range("L22").Value = 1773875
$ws.Range("M22").Value = 1819606
$ws.Range("D23").Value = 735
$ws.Range("E23").Value = 735
$ws.Range("F23").Value = 735
$ws.Range("G23").Value = 735
$ws.Range("H23").Value = 735
$ws.Range("I23").Value = 735
$ws.Range("J23").Value = 735
$ws.Range("K23").Value = 735
$ws.Range("L23").Value = 735
$ws.Range("M23").Value = 735
$ws.Range("D24").Value = '-'
$ws.Range("E24").Value = '-'
$ws.Range("F24").Value = '-'
$ws.Range("G24").Value = '-'
$ws.Range("H24").Value = '-'
$ws.Range("I24").Value = '-'
$ws.Range("J24").Value = '-'
$ws.Range("K24").Value = '-'
$ws.Range("L24").Value = '-'
$ws.Range("M24").Value = '-'
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 991840
$ws.Range("E26").Value = 1069072
$ws.Range("F26").Value = 1124152
$ws.Range("G26").Value = 1448077
$ws.Range("H26").Value = 2463375
$ws.Range("I26").Value = 2616297
$ws.Range("J26").Value = 2621368
$ws.Range("K26").Value = 2741472
$ws.Range("L26").Value = 2730371
$ws.Range("M26").Value = 2768140
$ws.Range("D27").Value = 7085298
$ws.Range("E27").Value = 8303575
$ws.Range("F27").Value = 7828612
$ws.Range("G27").Value = 9586494
$ws.Range("H27").Value = 10853828
$ws.Range("I27").Value = 12036422
$ws.Range("J27").Value = 12051163
$ws.Range("K27").Value = 12941826
$ws.Range("L27").Value = 13704436
$ws.Range("M27").Value = 17433361
$ws.Range("D29").Value = 659018
$ws.Range("E29").Value = 700865
$ws.Range("F29").Value = 407345
$ws.Range("G29").Value = 1092017
$ws.Range("H29").Value = 1648900
$ws.Range("I29").Value = 1576573
$ws.Range("J29").Value = 1454895
$ws.Range("K29").Value = 1350509
$ws.Range("L29").Value = 2082266
$ws.Range("M29").Value = 3776438
$ws.Range("D30").Value = '-'
$ws.Range("E30").Value = '-'
$ws.Range("F30").Value = '-'
$ws.Range("G30").Value = '-'
$ws.Range("H30").Value = '-'
$ws.Range("I30").Value = '-'
$ws.Range("J30").Value = '-'
$ws.Range("K30").Value = '-'
$ws.Range("L30").Value = '-'
$ws.Range("M30").Value = '-'
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = 149229
$ws.Range("E32").Value = 127844
$ws.Range("F32").Value = 127844
$ws.Range("G32").Value = 140688
$ws.Range("H32").Value = 232018
$ws.Range("I32").Value = 166358
$ws.Range("J32").Value = 167132
$ws.Range("K32").Value = 156411
$ws.Range("L32").Value = 170894
$ws.Range("M32").Value = 23713
$ws.Range("D33").Value = 971252
$ws.Range("E33").Value = 813079
$ws.Range("F33").Value = 1821476
$ws.Range("G33").Value = 1817192
$ws.Range("H33").Value = 1807723
$ws.Range("I33").Value = 1590981
$ws.Range("J33").Value = 2328317
$ws.Range("K33").Value = 2756720
$ws.Range("L33").Value = 2314886
$ws.Range("M33").Value = 1893451
$ws.Range("D34").Value = 1823318
$ws.Range("E34").Value = 2192487
$ws.Range("F34").Value = 1929021
$ws.Range("G34").Value = 2678694
$ws.Range("H34").Value = 2992867
$ws.Range("I34").Value = 2769535
$ws.Range("J34").Value = 3162685
$ws.Range("K34").Value = 2943634
$ws.Range("L34").Value = 3031600
$ws.Range("M34").Value = 3523641
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("D37").Value = 3602817
$ws.Range("E37").Value = 3834275
$ws.Range("F37").Value = 4285686
$ws.Range("G37").Value = 5728591
$ws.Range("H37").Value = 6681508
$ws.Range("I37").Value = 6103447
$ws.Range("J37").Value = 7113029
$ws.Range("K37").Value = 7207274
$ws.Range("L37").Value = 7599646
$ws.Range("M37").Value = 9217243
$ws.Range("D38").Value = 29550
$ws.Range("E38").Value = 34543
$ws.Range("F38").Value = 34543
$ws.Range("G38").Value = 41390
$ws.Range("H38").Value = 20453
$ws.Range("I38").Value = 48447
$ws.Range("J38").Value = 52673
$ws.Range("K38").Value = 58692
$ws.Range("L38").Value = 63656
$ws.Range("M38").Value = 93958
$ws.Range("D39").Value = '-'
$ws.Range("E39").Value = '-'
$ws.Range("F39").Value = '-'
$ws.Range("G39").Value = '-'
$ws.Range("H39").Value = '-'
$ws.Range("I39").Value = '-'
$ws.Range("J39").Value = '-'
$ws.Range("K39").Value = '-'
$ws.Range("L39").Value = '-'
$ws.Range("M39").Value = '-'
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("D41").Value = 80385
$ws.Range("E41").Value = 89632
$ws.Range("F41").Value = 126732
$ws.Range("G41").Value = 131138
$ws.Range("H41").Value = 134328
$ws.Range("I41").Value = 139815
$ws.Range("J41").Value = 218030
$ws.Range("K41").Value = 226210
$ws.Range("L41").Value = 230538
$ws.Range("M41").Value = 261212
$ws.Range("D42").Value = 109935
$ws.Range("E42").Value = 124175
$ws.Range("F42").Value = 161275
$ws.Range("G42").Value = 172528
$ws.Range("H42").Value = 154781
$ws.Range("I42").Value = 188262
$ws.Range("J42").Value = 270703
$ws.Range("K42").Value = 284902
$ws.Range("L42").Value = 294194
$ws.Range("M42").Value = 355170
$ws.Range("D43").Value = 3712752
$ws.Range("E43").Value = 3958450
$ws.Range("F43").Value = 4446961
$ws.Range("G43").Value = 5901119
$ws.Range("H43").Value = 6836289
$ws.Range("I43").Value = 6291709
$ws.Range("J43").Value = 7383732
$ws.Range("K43").Value = 7492176
$ws.Range("L43").Value = 7893840
$ws.Range("M43").Value = 9572413
$ws.Range("D45").Value = 1134000
$ws.Range("E45").Value = 2268000
$ws.Range("F45").Value = 2268000
$ws.Range("G45").Value = 2268000
$ws.Range("H45").Value = 2268000
$ws.Range("I45").Value = 2268000
$ws.Range("J45").Value = 2268000
$ws.Range("K45").Value = 2268000
$ws.Range("L45").Value = 2268000
$ws.Range("M45").Value = 2268000
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("D47").Value = 1134000
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("D50").Value = 113400
$ws.Range("E50").Value = 193764
$ws.Range("F50").Value = 193764
$ws.Range("G50").Value = 193764
$ws.Range("H50").Value = 193764
$ws.Range("I50").Value = 320443
$ws.Range("J50").Value = 320443
$ws.Range("K50").Value = 320443
$ws.Range("L50").Value = 320443
$ws.Range("M50").Value = 320443
$ws.Range("D51").Value = 97
$ws.Range("E51").Value = 97
$ws.Range("F51").Value = 97
$ws.Range("G51").Value = 97
$ws.Range("H51").Value = 97
$ws.Range("I51").Value = 97
$ws.Range("J51").Value = 97
$ws.Range("K51").Value = 97
$ws.Range("L51").Value = 97
$ws.Range("M51").Value = 97
$ws.Range("D52").Value = '-'
$ws.Range("E52").Value = '-'
$ws.Range("F52").Value = '-'
$ws.Range("G52").Value = '-'
$ws.Range("H52").Value = '-'
$ws.Range("I52").Value = '-'
$ws.Range("J52").Value = '-'
$ws.Range("K52").Value = '-'
$ws.Range("L52").Value = '-'
$ws.Range("M52").Value = '-'
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = '-'
$ws.Range("E54").Value = '-'
$ws.Range("F54").Value = '-'
$ws.Range("G54").Value = '-'
$ws.Range("H54").Value = '-'
$ws.Range("I54").Value = '-'
$ws.Range("J54").Value = '-'
$ws.Range("K54").Value = '-'
$ws.Range("L54").Value = '-'
$ws.Range("M54").Value = '-'
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("D56").Value = 991049
$ws.Range("E56").Value = 1883264
$ws.Range("F56").Value = 919790
$ws.Range("G56").Value = 1223514
$ws.Range("H56").Value = 1555678
$ws.Range("I56").Value = 3156173
$ws.Range("J56").Value = 2078891
$ws.Range("K56").Value = 2861110
$ws.Range("L56").Value = 3222056
$ws.Range("M56").Value = 5272408
$ws.Range("D57").Value = 3372546
$ws.Range("E57").Value = 4345125
$ws.Range("F57").Value = 3381651
$ws.Range("G57").Value = 3685375
$ws.Range("H57").Value = 4017539
$ws.Range("I57").Value = 5744713
$ws.Range("J57").Value = 4667431
$ws.Range("K57").Value = 5449650
$ws.Range("L57").Value = 5810596
$ws.Range("M57").Value = 7860948
$ws.Range("D58").Value = 7085298
$ws.Range("E58").Value = 8303575
$ws.Range("F58").Value = 7828612
$ws.Range("G58").Value = 9586494
$ws.Range("H58").Value = 10853828
$ws.Range("I58").Value = 12036422
$ws.Range("J58").Value = 12051163
$ws.Range("K58").Value = 12941826
$ws.Range("L58").Value = 13704436
$ws.Range("M58").Value = 17433361
